$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the price/volume columns as text so that numeric-looking
# strings (e.g. "1.005") are not auto-converted to numbers by Excel.
$fmtRange = $ws.Range("D2:E51")
$fmtRange.NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '26.428.14'
$ws.Cells.Item(2, 5).Value = '  +6.07%  '
$ws.Cells.Item(3, 4).Value = '1.719.77'
$ws.Cells.Item(3, 5).Value = '  +3.30%  '
$ws.Cells.Item(4, 4).Value = '1.005'
$ws.Cells.Item(4, 5).Value = '  +0.31%  '
$ws.Cells.Item(5, 4).Value = '330.95'
$ws.Cells.Item(5, 5).Value = '  +1.38%  '
$ws.Cells.Item(6, 4).Value = '1.002'
$ws.Cells.Item(6, 5).Value = '  +0.09%  '
$ws.Cells.Item(7, 4).Value = '0.3704'
$ws.Cells.Item(7, 5).Value = '  +1.97%  '
$ws.Cells.Item(8, 4).Value = '48.22'
$ws.Cells.Item(8, 5).Value = '  +1.05%  '
$ws.Cells.Item(9, 4).Value = '0.3347'
$ws.Cells.Item(9, 5).Value = '  +2.28%  '
$ws.Cells.Item(10, 5).Value = '  +3.78%  '
$ws.Cells.Item(11, 4).Value = '0.07358'
$ws.Cells.Item(11, 5).Value = '  +3.69%  '
$ws.Cells.Item(12, 4).Value = '1.001'
$ws.Cells.Item(12, 5).Value = '  +0.11%  '
$ws.Cells.Item(13, 4).Value = '6.369'
$ws.Cells.Item(13, 5).Value = '  +5.23%  '
$ws.Cells.Item(14, 4).Value = '20.02'
$ws.Cells.Item(14, 5).Value = '  +2.41%  '
$ws.Cells.Item(15, 4).Value = '7.008'
$ws.Cells.Item(15, 5).Value = '  +5.96%  '
$ws.Cells.Item(16, 4).Value = '1.722.83'
$ws.Cells.Item(16, 5).Value = '  +4.04%  '
$ws.Cells.Item(17, 4).Value = '0.00001067'
$ws.Cells.Item(17, 5).Value = '  +1.76%  '
$ws.Cells.Item(18, 4).Value = '0.06614'
$ws.Cells.Item(18, 5).Value = '  +0.08%  '
$ws.Cells.Item(19, 4).Value = '82.26'
$ws.Cells.Item(19, 5).Value = '  +3.94%  '
$ws.Cells.Item(20, 5).Value = '  +0.21%  '
$ws.Cells.Item(21, 5).Value = '  +4.39%  '
$ws.Cells.Item(22, 4).Value = '6.089'
$ws.Cells.Item(22, 5).Value = '  +2.75%  '
$ws.Cells.Item(23, 4).Value = '12.76'
$ws.Cells.Item(23, 5).Value = '  +1.30%  '
$ws.Cells.Item(24, 4).Value = '26.429.96'
$ws.Cells.Item(24, 5).Value = '  +6.25%  '
$ws.Cells.Item(25, 4).Value = '2.439'
$ws.Cells.Item(25, 5).Value = '  -0.48%  '
$ws.Cells.Item(26, 4).Value = '1.394'
$ws.Cells.Item(26, 5).Value = '  +17.48%  '
$ws.Cells.Item(27, 4).Value = '2.378'
$ws.Cells.Item(27, 5).Value = '  -2.36%  '
$ws.Cells.Item(28, 4).Value = '151.92'
$ws.Cells.Item(28, 5).Value = '  +2.02%  '
$ws.Cells.Item(29, 4).Value = '19.34'
$ws.Cells.Item(29, 5).Value = '  +3.59%  '
$ws.Cells.Item(30, 4).Value = '1.916.85'
$ws.Cells.Item(30, 5).Value = '  +4.27%  '
$ws.Cells.Item(31, 4).Value = '130.48'
$ws.Cells.Item(31, 5).Value = '  +3.89%  '
$ws.Cells.Item(32, 4).Value = '4.105'
$ws.Cells.Item(32, 5).Value = '  +0.39%  '
$ws.Cells.Item(33, 4).Value = '5.939'
$ws.Cells.Item(33, 5).Value = '  +3.34%  '
$ws.Cells.Item(34, 4).Value = '0.08559'
$ws.Cells.Item(34, 5).Value = '  +1.36%  '
$ws.Cells.Item(35, 4).Value = '1.695'
$ws.Cells.Item(35, 5).Value = '  +2.70%  '
$ws.Cells.Item(36, 4).Value = '12.63'
$ws.Cells.Item(36, 5).Value = '  +3.28%  '
$ws.Cells.Item(37, 4).Value = '5.321'
$ws.Cells.Item(37, 5).Value = '  +2.64%  '
$ws.Cells.Item(38, 4).Value = '0.02309'
$ws.Cells.Item(38, 5).Value = '  +1.56%  '
$ws.Cells.Item(39, 5).Value = '  +3.46%  '
$ws.Cells.Item(40, 4).Value = '0.06174'
$ws.Cells.Item(40, 5).Value = '  -0.04%  '
$ws.Cells.Item(41, 4).Value = '8.419'
$ws.Cells.Item(41, 5).Value = '  +1.27%  '
$ws.Cells.Item(42, 4).Value = '1.222'
$ws.Cells.Item(42, 5).Value = '  -4.60%  '
$ws.Cells.Item(43, 4).Value = '0.6159'
$ws.Cells.Item(43, 5).Value = '  +3.66%  '
$ws.Cells.Item(44, 4).Value = '1.001'
$ws.Cells.Item(44, 5).Value = '  +0.12%  '
$ws.Cells.Item(45, 4).Value = '14.01'
$ws.Cells.Item(45, 5).Value = '  +3.63%  '
$ws.Cells.Item(46, 4).Value = '3.899'
$ws.Cells.Item(46, 5).Value = '  +1.65%  '
$ws.Cells.Item(47, 4).Value = '0.5941'
$ws.Cells.Item(47, 5).Value = '  +5.43%  '
$ws.Cells.Item(48, 4).Value = '127.50'
$ws.Cells.Item(48, 5).Value = '  +1.81%  '
$ws.Cells.Item(49, 4).Value = '2.027'
$ws.Cells.Item(49, 5).Value = '  +3.88%  '
$ws.Cells.Item(50, 4).Value = '0.07167'
$ws.Cells.Item(50, 5).Value = '  +2.52%  '
$ws.Cells.Item(51, 4).Value = '76.44'
$ws.Cells.Item(51, 5).Value = '  +1.97%  '

# Remove the temporary text formatting so cell styles match the original
# (unstyled) cells exactly.
$fmtRange.ClearFormats()
